$d = $word.ActiveDocument

# 1. "Recently, philosopher Kevin Vallier published an essay" -> add "political"
$d.Content.Find.Execute(
    "Recently, philosopher Kevin Vallier published an essay",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Recently, political philosopher Kevin Vallier published an essay", 2)

# 2. "which argues that," -> "in which he argues that,"
$d.Content.Find.Execute(
    "which argues that,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "in which he argues that,", 2)

# 3. "blue laws legislating strict limits to what commerce could" -> "...limits as to what..."
$d.Content.Find.Execute(
    "19th-century that blue laws legislating strict limits to what commerce could",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "19th-century that blue laws legislating strict limits as to what commerce could", 2)

# 4. "if one would not work on Sundays." -> "if one will not work on Sundays."
$d.Content.Find.Execute(
    "if one would not work on Sundays. (Max Weber described this as the",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "if one will not work on Sundays. (Max Weber described this as the", 2)

# 5a. "little shelters that express for their inhabitants how they" -> add comma after "shelters"
$d.Content.Find.Execute(
    "little shelters that express for their inhabitants how they",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "little shelters, that express for their inhabitants how they", 2)

# 5b. "are connected to and take part in," -> add comma after "to"
$d.Content.Find.Execute(
    "are connected to and take part in,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "are connected to, and take part in,", 2)

# 6. "mysteriously, inhabiting." -> "willy-nilly, inhabiting."
$d.Content.Find.Execute(
    "mysteriously, inhabiting.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "willy-nilly, inhabiting.", 2)
